$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 44.25
$ws.Range("I9").Value = 15
$ws.Range("J9").Value = 73.5
$ws.Range("K9").Value = 15
$ws.Range("L9").Value = 73.5
$ws.Range("M9").Value = 154
$ws.Range("N9").Value = -411.5
$ws.Range("H116").Value = 2330.3
$ws.Range("I116").Value = 2142.8572
$ws.Range("J116").Value = 2767.6667
$ws.Range("K116").Value = 2142.8572
$ws.Range("L116").Value = 2767.6667
$ws.Range("M116").Value = 1299.1428
$ws.Range("N116").Value = -9651.6667
$ws.Range("H138").Value = 5468256.5
$ws.Range("I138").Value = 10418926
$ws.Range("J138").Value = 5448.517
$ws.Range("K138").Value = 31256778
$ws.Range("L138").Value = 16345.551
$ws.Range("M138").Value = -31251638
$ws.Range("N138").Value = -26625.551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16208.138
$ws.Range("I32").Value = 15758.216
$ws.Range("J32").Value = 18769.23
$ws.Range("K32").Value = 15758.216
$ws.Range("L32").Value = 18769.23
$ws.Range("M32").Value = -15471.216
$ws.Range("N32").Value = -19343.23
$ws.Range("H45").Value = 1895509.1
$ws.Range("I45").Value = 2526828.8
$ws.Range("J45").Value = 1550
$ws.Range("K45").Value = 2526828.8
$ws.Range("L45").Value = 1550
$ws.Range("M45").Value = -2526451.8
$ws.Range("N45").Value = -2304
$ws.Range("H74").Value = 1376.421
$ws.Range("I74").Value = 1352.8889
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 1352.8889
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -478.8888999999999
$ws.Range("N74").Value = -3548
$ws.Range("H77").Value = 1376.421
$ws.Range("I77").Value = 1352.8889
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 6764.4445
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -2396.4445
$ws.Range("N77").Value = -17736
$ws.Range("H132").Value = 1863.7693
$ws.Range("I132").Value = 1425.6389
$ws.Range("J132").Value = 2849.5625
$ws.Range("K132").Value = 4276.9167
$ws.Range("L132").Value = 8548.6875
$ws.Range("M132").Value = -1746.9167
$ws.Range("N132").Value = -13608.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 86990
$ws.Range("J59").Value = 86990
$ws.Range("L59").Value = 86990
$ws.Range("N59").Value = -88684
$ws.Range("H118").Value = 7889.4736
$ws.Range("J118").Value = 7889.4736
$ws.Range("L118").Value = 7889.4736
$ws.Range("N118").Value = -11203.4736
$ws.Range("H132").Value = 37490.75
$ws.Range("J132").Value = 37490.75
$ws.Range("L132").Value = 37490.75
$ws.Range("N132").Value = -47610.75
$ws.Range("H134").Value = 2454.9214
$ws.Range("I134").Value = 1564.4918
$ws.Range("J134").Value = 4394.7856
$ws.Range("K134").Value = 4693.4754
$ws.Range("L134").Value = 13184.3568
$ws.Range("M134").Value = -2158.4754
$ws.Range("N134").Value = -18254.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1833.7059
$ws.Range("I122").Value = 2099.0908
$ws.Range("J122").Value = 1347.1666
$ws.Range("K122").Value = 6297.2724
$ws.Range("L122").Value = 4041.4998
$ws.Range("M122").Value = -3847.2724
$ws.Range("N122").Value = -8941.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 824.375
$ws.Range("I5").Value = 321.66666
$ws.Range("J5").Value = 1470.7142
$ws.Range("K5").Value = 964.9999799999999
$ws.Range("L5").Value = 4412.142599999999
$ws.Range("M5").Value = -852.9999799999999
$ws.Range("N5").Value = -4636.142599999999
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H107").Value = 547.4
$ws.Range("J107").Value = 646.0714
$ws.Range("L107").Value = 1938.2142
$ws.Range("N107").Value = -5778.2142
$ws.Range("H131").Value = 2824.1614
$ws.Range("I131").Value = 3329.9333
$ws.Range("J131").Value = 2350
$ws.Range("K131").Value = 9989.7999
$ws.Range("L131").Value = 7050
$ws.Range("M131").Value = -4949.7999
$ws.Range("N131").Value = -17130
$ws.Range("H135").Value = 824.375
$ws.Range("I135").Value = 321.66666
$ws.Range("J135").Value = 1470.7142
$ws.Range("K135").Value = 2894.99994
$ws.Range("L135").Value = 13236.4278
$ws.Range("M135").Value = -359.9999399999997
$ws.Range("N135").Value = -18306.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6900
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 6900
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 6900
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -7480
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("H126").Value = 4307.2856
$ws.Range("I126").Value = 3002.4
$ws.Range("K126").Value = 9007.200000000001
$ws.Range("M126").Value = -6537.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8226.6
$ws.Range("I7").Value = 6222.25
$ws.Range("J7").Value = 10076.77
$ws.Range("K7").Value = 6222.25
$ws.Range("L7").Value = 10076.77
$ws.Range("M7").Value = -6110.25
$ws.Range("N7").Value = -10300.77
$ws.Range("H33").Value = 7887.5
$ws.Range("I33").Value = 1300
$ws.Range("J33").Value = 10083.333
$ws.Range("K33").Value = 1300
$ws.Range("L33").Value = 10083.333
$ws.Range("M33").Value = -1010
$ws.Range("N33").Value = -10663.333
$ws.Range("H40").Value = 2914.9473
$ws.Range("I40").Value = 6025
$ws.Range("J40").Value = 1647.8889
$ws.Range("K40").Value = 6025
$ws.Range("L40").Value = 1647.8889
$ws.Range("M40").Value = -5889
$ws.Range("N40").Value = -1919.8889
$ws.Range("H126").Value = 8226.6
$ws.Range("I126").Value = 6222.25
$ws.Range("J126").Value = 10076.77
$ws.Range("K126").Value = 18666.75
$ws.Range("L126").Value = 30230.31
$ws.Range("M126").Value = -16196.75
$ws.Range("N126").Value = -35170.31
$ws.Range("H132").Value = 11912498
$ws.Range("I132").Value = 6070.2964
$ws.Range("J132").Value = 33344066
$ws.Range("K132").Value = 18210.8892
$ws.Range("L132").Value = 100032198
$ws.Range("M132").Value = -15680.8892
$ws.Range("N132").Value = -100037258

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 80811
$ws.Range("J108").Value = 80811
$ws.Range("L108").Value = 80811
$ws.Range("N108").Value = -88491
$ws.Range("H126").Value = 1717.0889
$ws.Range("I126").Value = 1850.56
$ws.Range("J126").Value = 1550.25
$ws.Range("K126").Value = 5551.68
$ws.Range("L126").Value = 4650.75
$ws.Range("M126").Value = -3081.68
$ws.Range("N126").Value = -9590.75
$ws.Range("H132").Value = 1867.7091
$ws.Range("I132").Value = 1696.6666
$ws.Range("J132").Value = 2420.3076
$ws.Range("K132").Value = 5089.9998
$ws.Range("L132").Value = 7260.9228
$ws.Range("M132").Value = -2559.9998
$ws.Range("N132").Value = -12320.9228
$ws.Range("H136").Value = 1325.277
$ws.Range("I136").Value = 640.39655
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 1921.18965
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = 628.8103499999997
$ws.Range("N136").Value = -26100
$ws.Range("H138").Value = 50096.89
$ws.Range("J138").Value = 50096.89
$ws.Range("L138").Value = 50096.89
$ws.Range("N138").Value = -60376.89
